# --- rial_cumulative.xlsx: roll the quarterly series forward one
# column (drop the oldest period, append the newly published one)
# and correct the re-published date for the 6m/1401-06 report. ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9
$ws.Range("D9").Value = "1400-10-30 (2)"
$ws.Range("E9").Value = "1401-05-19 (8)"
$ws.Range("F9").Value = "1401-04-29 (3)"
$ws.Range("G9").Value = "1401-09-29 (4)"
$ws.Range("H9").Value = "1401-11-01 (2)"
$ws.Range("I9").Value = "1402-02-27 (7)"
$ws.Range("K9").Value = "1401-09-29 (2)"

# Row 11
$ws.Range("D11").Value = 3352884
$ws.Range("E11").Value = 5321509
$ws.Range("F11").Value = 2176989
$ws.Range("G11").Value = 4690774
$ws.Range("H11").Value = 7636647
$ws.Range("I11").Value = 10593475
$ws.Range("J11").Value = 5670880
$ws.Range("K11").Value = 10287179
$ws.Range("L11").Value = 14022766
$ws.Range("M11").Value = 20405785

# Row 12
$ws.Range("D12").Value = -2490143
$ws.Range("E12").Value = -3914091
$ws.Range("F12").Value = -1471721
$ws.Range("G12").Value = -3366689
$ws.Range("H12").Value = -5761186
$ws.Range("I12").Value = -8499119
$ws.Range("J12").Value = -4545645
$ws.Range("K12").Value = -7896360
$ws.Range("L12").Value = -10755266
$ws.Range("M12").Value = -16273879

# Row 13
$ws.Range("D13").Value = 862741
$ws.Range("E13").Value = 1407418
$ws.Range("F13").Value = 705268
$ws.Range("G13").Value = 1324085
$ws.Range("H13").Value = 1875461
$ws.Range("I13").Value = 2094356
$ws.Range("J13").Value = 1125235
$ws.Range("K13").Value = 2390819
$ws.Range("L13").Value = 3267500
$ws.Range("M13").Value = 4131906

# Row 14
$ws.Range("D14").Value = -181026
$ws.Range("E14").Value = -264863
$ws.Range("F14").Value = -61623
$ws.Range("G14").Value = -153543
$ws.Range("H14").Value = -231837
$ws.Range("I14").Value = -381452
$ws.Range("J14").Value = -172347
$ws.Range("K14").Value = -331665
$ws.Range("L14").Value = -522814
$ws.Range("M14").Value = -984898

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = 12127
$ws.Range("E16").Value = 13693
$ws.Range("F16").Value = 3103
$ws.Range("G16").Value = 13999
$ws.Range("H16").Value = -5539
$ws.Range("I16").Value = 23054
$ws.Range("J16").Value = -10029
$ws.Range("K16").Value = 14091
$ws.Range("L16").Value = -20103
$ws.Range("M16").Value = 20262

# Row 17
$ws.Range("D17").Value = 693842
$ws.Range("E17").Value = 1156248
$ws.Range("F17").Value = 646748
$ws.Range("G17").Value = 1184541
$ws.Range("H17").Value = 1638085
$ws.Range("I17").Value = 1735958
$ws.Range("J17").Value = 942859
$ws.Range("K17").Value = 2073245
$ws.Range("L17").Value = 2724583
$ws.Range("M17").Value = 3167270

# Row 18
$ws.Range("D18").Value = -16474
$ws.Range("E18").Value = -20691
$ws.Range("F18").Value = -5169
$ws.Range("G18").Value = -12054
$ws.Range("H18").Value = -17070
$ws.Range("I18").Value = -39361
$ws.Range("J18").Value = -54830
$ws.Range("K18").Value = -220920
$ws.Range("L18").Value = -361711
$ws.Range("M18").Value = -435585

# Row 19
$ws.Range("D19").Value = 4937
$ws.Range("E19").Value = -458290
$ws.Range("F19").Value = 1152
$ws.Range("G19").Value = 1001
$ws.Range("H19").Value = 3719
$ws.Range("I19").Value = 499043
$ws.Range("J19").Value = -1429
$ws.Range("K19").Value = 2627
$ws.Range("L19").Value = -331
$ws.Range("M19").Value = -55389

# Row 20
$ws.Range("D20").Value = 682305
$ws.Range("E20").Value = 677267
$ws.Range("F20").Value = 642731
$ws.Range("G20").Value = 1173488
$ws.Range("H20").Value = 1624734
$ws.Range("I20").Value = 2195640
$ws.Range("J20").Value = 886600
$ws.Range("K20").Value = 1854952
$ws.Range("L20").Value = 2362541
$ws.Range("M20").Value = 2676296

# Row 21
$ws.Range("D21").Value = -21
$ws.Range("E21").Value = -30
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("D22").Value = 682284
$ws.Range("E22").Value = 677237
$ws.Range("F22").Value = 642731
$ws.Range("G22").Value = 1173488
$ws.Range("H22").Value = 1624734
$ws.Range("I22").Value = 2195640
$ws.Range("J22").Value = 886600
$ws.Range("K22").Value = 1854952
$ws.Range("L22").Value = 2362541
$ws.Range("M22").Value = 2676296

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 682284
$ws.Range("E24").Value = 677237
$ws.Range("F24").Value = 642731
$ws.Range("G24").Value = 1173488
$ws.Range("H24").Value = 1624734
$ws.Range("I24").Value = 2195640
$ws.Range("J24").Value = 886600
$ws.Range("K24").Value = 1854952
$ws.Range("L24").Value = 2362541
$ws.Range("M24").Value = 2676296

# Row 25
$ws.Range("D25").Value = 975
$ws.Range("E25").Value = 967
$ws.Range("F25").Value = 918
$ws.Range("G25").Value = 1676
$ws.Range("H25").Value = 2321
$ws.Range("I25").Value = 3137
$ws.Range("J25").Value = 1267
$ws.Range("K25").Value = 2650
$ws.Range("L25").Value = 3375
$ws.Range("M25").Value = 3823

# Row 26
$ws.Range("D26").Value = 700000
$ws.Range("E26").Value = 700000
$ws.Range("F26").Value = 700000
$ws.Range("G26").Value = 700000
$ws.Range("H26").Value = 700000
$ws.Range("I26").Value = 700000
$ws.Range("J26").Value = 700000
$ws.Range("K26").Value = 700000
$ws.Range("L26").Value = 700000
$ws.Range("M26").Value = 700000

# Row 27
$ws.Range("D27").Value = 975
$ws.Range("E27").Value = 967
$ws.Range("F27").Value = 918
$ws.Range("G27").Value = 1676
$ws.Range("H27").Value = 2321
$ws.Range("I27").Value = 3137
$ws.Range("J27").Value = 1267
$ws.Range("K27").Value = 2650
$ws.Range("L27").Value = 3375
$ws.Range("M27").Value = 3823

# A few cells hold a bare ISO date ("YYYY-MM-DD"); written as plain
# text they would otherwise be auto-parsed into a date serial, so
# force Text format first, then restore the normal cell style by
# pasting formats only from an already-correct sibling cell.
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1401-11-01"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-27"
$ws.Range("K9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("M9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths also roll forward one column (D:M carry the 10
# quarterly series; widths alternate 28/29 every 4th column).
$ws.Columns("D:D").ColumnWidth = 27.1
$ws.Columns("E:E").ColumnWidth = 28.1
$ws.Columns("F:H").ColumnWidth = 27.1
$ws.Columns("I:I").ColumnWidth = 28.1
$ws.Columns("J:L").ColumnWidth = 27.1
$ws.Columns("M:M").ColumnWidth = 28.1